$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "LogicalName"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Id"
$ws.Range("D1").Value = "Xpath"
$ws.Range("E1").Value = "Css"

# Fill column by column so the shared-string table is built in the
# same order as the source workbook (headers, then col A, col B, ...)
$ws.Range("A2").Value = "obj1"
$ws.Range("A3").Value = "obj2"
$ws.Range("A4").Value = "obj3"

$ws.Range("B2").Value = "name1"
$ws.Range("B3").Value = "name2"
$ws.Range("B4").Value = "name2"

$ws.Range("C2").Value = "id1"
$ws.Range("C3").Value = "id2"
$ws.Range("C4").Value = "id3"

$ws.Range("D2").Value = "xpath1"
$ws.Range("D3").Value = "xpath2"
$ws.Range("D4").Value = "xpath3"

$ws.Range("E2").Value = "css1"
$ws.Range("E3").Value = "css2"
$ws.Range("E4").Value = "css3"

$ws.Range("E4").Select()
